$d = $word.ActiveDocument

# Set single line spacing (wdLineSpaceSingle = 0) for every paragraph in
# the document body. This writes <w:spacing w:line="240" w:lineRule="auto"/>
# into each paragraph's pPr (creating pPr where needed), matching Word's
# "Line Spacing = Single" formatting.
$d.Paragraphs.LineSpacingRule = 0
